# Edit 1: "Males are more likely to be in a romantic relationship" ->
#         "Females are more likely to be in a romantic relationship"
# The target XML keeps this split across two runs: "Fem" + "ales are more
# likely to be in a romantic relationship", with the trailing "." staying
# in its own (pre-existing) run.

$d = $word.ActiveDocument

$full = $d.Content
$found = $full.Find.Execute("Males are more likely to be in a romantic relationship", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $full.Start
$end = $full.End

# Toggle Bold briefly on the boundary pieces so the engine's run-coalescing
# pass (which fires on any Delete) does not re-merge the runs we are about
# to split; we flip it back to its original value immediately after.

# Protect the boundary right after this sentence (the following "." run).
$dotRange = $d.Range($end, $end + 1)
$dotRange.Bold = 1

# Split off "Fem" as its own leading run.
$insPoint = $d.Range($start, $start)
$insPoint.InsertBefore("Fem")

$femRange = $d.Range($start, $start + 3)
$femRange.Bold = 1

# Remove the "M" of "Males" so the second run reads "ales are more likely...".
$mRange = $d.Range($start + 3, $start + 4)
$mRange.Delete()

# Restore formatting on both protected pieces.
$femRange2 = $d.Range($start, $start + 3)
$femRange2.Bold = 0

$dotRange2 = $d.Range($end + 3 - 1, $end + 3 - 1 + 1)
$dotRange2.Bold = 0

# Edit 2: move the (hidden) "_GoBack" bookmark from between "mothe" / "r's
# education level." (next bullet) to between "Positive relationship
# between " / "G3 scores and mother and father education levels." (this
# bullet) -- a side effect of where the author's cursor was after the last
# edit. Re-adding a bookmark with the same name moves it (removing the old
# one).

$p1 = $d.Content
$m1 = $p1.Find.Execute("Positive relationship between ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
$p2 = $d.Range($p1.End, $d.Content.End)
$m2 = $p2.Find.Execute("Positive relationship between ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

$bmPoint = $d.Range($p2.End, $p2.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
